$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.712771666666666
$ws.Range("H2").Value = 17.138315
$ws.Range("I2").Value = 0.1683613830606884
$ws.Range("J2").Value = 0.1683613830606885
$ws.Range("M2").Value = 13.35941066666667
$ws.Range("N2").Value = 40.078232
$ws.Range("O2").Value = 0.4925555025958562
$ws.Range("P2").Value = 0.4925555025958562
$ws.Range("Q2").Value = 76.31926273989777
$ws.Range("R2").Value = 686.8733646590799
$ws.Range("S2").Value = 0.08292732565119086
$ws.Range("T2").Value = 0.08292732565119088
$ws.Range("G3").Value = 5.712771666666666
$ws.Range("H3").Value = 17.138315
$ws.Range("I3").Value = 0.1683613830606884
$ws.Range("J3").Value = 0.1683613830606885
$ws.Range("O3").Value = 0.03774352140193379
$ws.Range("P3").Value = 0.03774352140193379
$ws.Range("Q3").Value = 5.848189110510555
$ws.Range("R3").Value = 52.63370199459499
$ws.Range("S3").Value = 0.006354551464810267
$ws.Range("T3").Value = 0.006354551464810267
$ws.Range("G4").Value = 5.712771666666666
$ws.Range("H4").Value = 17.138315
$ws.Range("I4").Value = 0.1683613830606884
$ws.Range("J4").Value = 0.1683613830606885
$ws.Range("M4").Value = 12.73953533333333
$ws.Range("N4").Value = 38.218606
$ws.Range("O4").Value = 0.4697009760022101
$ws.Range("P4").Value = 0.46970097600221
$ws.Range("Q4").Value = 72.77805649876555
$ws.Range("R4").Value = 655.0025084888899
$ws.Range("S4").Value = 0.07907950594468731
$ws.Range("T4").Value = 0.07907950594468731
$ws.Range("I5").Value = 0.4370667227533506
$ws.Range("J5").Value = 0.4370667227533506
$ws.Range("M5").Value = 13.35941066666667
$ws.Range("N5").Value = 40.078232
$ws.Range("O5").Value = 0.4925555025958562
$ws.Range("P5").Value = 0.4925555025958562
$ws.Range("Q5").Value = 198.1250655125298
$ws.Range("R5").Value = 1783.125589612768
$ws.Range("S5").Value = 0.2152796192937003
$ws.Range("T5").Value = 0.2152796192937003
$ws.Range("I6").Value = 0.4370667227533506
$ws.Range("J6").Value = 0.4370667227533506
$ws.Range("O6").Value = 0.03774352140193379
$ws.Range("P6").Value = 0.03774352140193379
$ws.Range("S6").Value = 0.01649643720431415
$ws.Range("T6").Value = 0.01649643720431415
$ws.Range("I7").Value = 0.4370667227533506
$ws.Range("J7").Value = 0.4370667227533506
$ws.Range("M7").Value = 12.73953533333333
$ws.Range("N7").Value = 38.218606
$ws.Range("O7").Value = 0.4697009760022101
$ws.Range("P7").Value = 0.46970097600221
$ws.Range("Q7").Value = 188.9320820725715
$ws.Range("R7").Value = 1700.388738653144
$ws.Range("S7").Value = 0.2052906662553361
$ws.Range("T7").Value = 0.2052906662553361
$ws.Range("G8").Value = 13.388457
$ws.Range("H8").Value = 40.165371
$ws.Range("I8").Value = 0.3945718941859609
$ws.Range("J8").Value = 0.3945718941859609
$ws.Range("M8").Value = 13.35941066666667
$ws.Range("N8").Value = 40.078232
$ws.Range("O8").Value = 0.4925555025958562
$ws.Range("P8").Value = 0.4925555025958562
$ws.Range("Q8").Value = 178.861895256008
$ws.Range("R8").Value = 1609.757057304072
$ws.Range("S8").Value = 0.194348557650965
$ws.Range("T8").Value = 0.194348557650965
$ws.Range("G9").Value = 13.388457
$ws.Range("H9").Value = 40.165371
$ws.Range("I9").Value = 0.3945718941859609
$ws.Range("J9").Value = 0.3945718941859609
$ws.Range("O9").Value = 0.03774352140193379
$ws.Range("P9").Value = 0.03774352140193379
$ws.Range("Q9").Value = 13.705821447547
$ws.Range("R9").Value = 123.352393027923
$ws.Range("S9").Value = 0.01489253273280937
$ws.Range("T9").Value = 0.01489253273280937
$ws.Range("G10").Value = 13.388457
$ws.Range("H10").Value = 40.165371
$ws.Range("I10").Value = 0.3945718941859609
$ws.Range("J10").Value = 0.3945718941859609
$ws.Range("M10").Value = 12.73953533333333
$ws.Range("N10").Value = 38.218606
$ws.Range("O10").Value = 0.4697009760022101
$ws.Range("P10").Value = 0.46970097600221
$ws.Range("Q10").Value = 170.562721010314
$ws.Range("R10").Value = 1535.064489092826
$ws.Range("S10").Value = 0.1853308038021866
$ws.Range("T10").Value = 0.1853308038021866
